$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("201:201").Insert()

$ws.Range("A201").Value = 5
$ws.Range("B201").Value = "Macroferia Regional de Talca"
$ws.Range("C201").Value = "Maule"
$ws.Range("D201").Value = 44609
$ws.Range("E201").Value = 7
$ws.Range("F201").Value = 100112032
$ws.Range("G201").Value = "Zapallo italiano"
$ws.Range("H201").Value = "Sin especificar"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 400
$ws.Range("K201").Value = 6000
$ws.Range("L201").Value = 6000
$ws.Range("M201").Value = 6000
$ws.Range("N201").Value = "`$/caja 50 unidades"
$ws.Range("O201").Value = "Región del Maule"
$ws.Range("P201").Value = 120
$ws.Range("Q201").Value = 50
$ws.Range("R201").Value = "Hortaliza"
